# code and data cleanup
#  - rename EMPLOYMENT_DATA -> UNEMPLOYMENT_DATA (and its header/shared string)
#  - move the active selection around (GDP_DATA becomes the active tab,
#    selection on the renamed sheet moves to F11)

$wb = $excel.ActiveWorkbook

# Rename the sheet and fix the column header text that referenced the old name.
$wsUnemp = $wb.Worksheets.Item("EMPLOYMENT_DATA")
$wsUnemp.Name = "UNEMPLOYMENT_DATA"
$wsUnemp.Range("E1").Value = "UNEMPLOYMENT_RATE (Percent)"

# Update the selection left on the renamed sheet.
$wsUnemp.Range("F11").Select() | Out-Null

# Make GDP_DATA the active sheet/tab with cell A1 selected.
$wsGdp = $wb.Worksheets.Item("GDP_DATA")
$wsGdp.Activate() | Out-Null
$wsGdp.Range("A1").Select() | Out-Null
